# Cap nhat 10 bai toan thu nghiem
# Update experiment result values on Sheet1 to reflect the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Problem 2 (row 5)
$ws.Range("B5").Value = 49
$ws.Range("C5").Value = 10.44
$ws.Range("D5").Value = 30
$ws.Range("F5").Value = 75
$ws.Range("L5").Value = 75.5
$ws.Range("M5").Value = 5.22

# Problem 3 (row 6)
$ws.Range("B6").Value = 19.5
$ws.Range("C6").Value = 5.68
$ws.Range("F6").Value = 30

# Problem 4 (row 7)
$ws.Range("B7").Value = 31.8
$ws.Range("C7").Value = 1.99

# Problem 5 (row 8)
$ws.Range("B8").Value = 123.9
$ws.Range("C8").Value = 15.69
$ws.Range("D8").Value = 93
$ws.Range("Q8").Value = 239

# Problem 6 (row 9)
$ws.Range("B9").Value = 111.6
$ws.Range("C9").Value = 5.54
$ws.Range("L9").Value = 121.1
$ws.Range("M9").Value = 6.46
$ws.Range("N9").Value = 110

# Problem 7 (row 10)
$ws.Range("B10").Value = 105.8
$ws.Range("C10").Value = 5.27
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = 102
$ws.Range("L10").Value = 114.6
$ws.Range("M10").Value = 1.96
$ws.Range("N10").Value = 109

# Problem 8 (row 11)
$ws.Range("B11").Value = 54.3
$ws.Range("C11").Value = 4.78
$ws.Range("D11").Value = 49
$ws.Range("E11").Value = 52
$ws.Range("F11").Value = 61
$ws.Range("G11").Value = 23
$ws.Range("I11").Value = 23
$ws.Range("J11").Value = 23
$ws.Range("K11").Value = 23
$ws.Range("L11").Value = 84.5
$ws.Range("M11").Value = 15.5
$ws.Range("N11").Value = 69
$ws.Range("O11").Value = 84.5

# Problem 9 (row 12)
$ws.Range("B12").Value = 111.1
$ws.Range("C12").Value = 17.55
$ws.Range("D12").Value = 98
$ws.Range("E12").Value = 107
$ws.Range("F12").Value = 162
$ws.Range("G12").Value = 130
$ws.Range("I12").Value = 130
$ws.Range("J12").Value = 130
$ws.Range("K12").Value = 130
$ws.Range("L12").Value = 174.3
$ws.Range("M12").Value = 4.2
$ws.Range("N12").Value = 163
$ws.Range("O12").Value = 176
$ws.Range("P12").Value = 177

# Problem 10 (row 13)
$ws.Range("B13").Value = 132
$ws.Range("C13").Value = 6
$ws.Range("E13").Value = 130
$ws.Range("L13").Value = 142
$ws.Range("M13").Value = 6
$ws.Range("N13").Value = 130
$ws.Range("P13").Value = 150
